$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8166753333333334
$ws.Range("H2").Value = 2.450026
$ws.Range("M2").Value = 30.224788
$ws.Range("N2").Value = 90.674364
$ws.Range("O2").Value = 0.3247052378228209
$ws.Range("P2").Value = 0.3247052378228209
$ws.Range("Q2").Value = 24.68383881482934
$ws.Range("R2").Value = 222.154549333464
$ws.Range("S2").Value = 0.3247052378228209
$ws.Range("T2").Value = 0.3247052378228209

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8166753333333334
$ws.Range("H3").Value = 2.450026
$ws.Range("M3").Value = 20.25845733333333
$ws.Range("N3").Value = 60.775372
$ws.Range("O3").Value = 0.2176368352473959
$ws.Range("P3").Value = 0.217636835247396
$ws.Range("Q3").Value = 16.54458239551911
$ws.Range("R3").Value = 148.901241559672
$ws.Range("S3").Value = 0.2176368352473959
$ws.Range("T3").Value = 0.217636835247396

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8166753333333334
$ws.Range("H4").Value = 2.450026
$ws.Range("M4").Value = 12.725178
$ws.Range("N4").Value = 38.175534
$ws.Range("O4").Value = 0.1367067305427495
$ws.Range("P4").Value = 0.1367067305427495
$ws.Range("Q4").Value = 10.392338984876
$ws.Range("R4").Value = 93.531050863884
$ws.Range("S4").Value = 0.1367067305427495
$ws.Range("T4").Value = 0.1367067305427495

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8166753333333334
$ws.Range("H5").Value = 2.450026
$ws.Range("M5").Value = 29.87534766666667
$ws.Range("N5").Value = 89.62604300000001
$ws.Range("O5").Value = 0.3209511963870337
$ws.Range("P5").Value = 0.3209511963870337
$ws.Range("Q5").Value = 24.39845951412423
$ws.Range("R5").Value = 219.586135627118
$ws.Range("S5").Value = 0.3209511963870337
$ws.Range("T5").Value = 0.3209511963870337
